$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) and Volume(1h) (E) columns are refreshed with the latest
# crypto market snapshot. D-column values are plain numeric-looking
# strings (e.g. "1.00", "36.00") that must stay as literal text, so we
# force a text number format before assigning them, then restore the
# default "Normal" style so no stray style index is left behind.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.127.18"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.16%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.312.96"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.84%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.15%  "
$ws.Range("E7").Value = "  +2.37%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.525"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.87%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.00"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0817"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.75%  "
$ws.Range("E12").Value = "  +1.01%  "
$ws.Range("E13").Value = "  +6.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.673.37"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.98"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.59%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.321.22"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.47%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.812"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.075.33"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.57"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.43%  "
$ws.Range("E20").Value = "  +2.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.13"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.50"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.70%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "241.23"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.63"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.84"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.71%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.39"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.69%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.68"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.21%  "
$ws.Range("E30").Value = "  +0.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "166.49"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.34"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.51%  "
$ws.Range("E33").Value = "  +0.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.15"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "18.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.72%  "
$ws.Range("E36").Value = "  +1.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.107"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.56%  "
$ws.Range("E39").Value = "  +3.03%  "
$ws.Range("E40").Value = "  +2.44%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.30"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.70%  "
$ws.Range("E42").Value = "  +1.66%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "19.58"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0290"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.979.87"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.01"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.84"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.96"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +18.80%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "55.58"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.77%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.543.54"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.96%  "
$ws.Range("E51").Value = "  +4.08%  "
